$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$players = @(
    "Malik Monk",
    "Cade Cunningham",
    "Derrick White",
    "Herbert Jones",
    "Duncan Robinson",
    "Kelly Olynyk",
    "Bam Adebayo",
    "Isaiah Hartenstein",
    "Anthony Davis",
    "Julius Randle",
    "Damian Lillard",
    "Kentavious Caldwell-Pope",
    "Cameron Johnson",
    "Bilal Coulibaly",
    "Brandon Miller",
    "Brandon Ingram",
    "LaMelo Ball"
)

$positions = @(
    "SG,SF",
    "PG,SG",
    "PG,SG",
    "SF,PF",
    "SG,SF",
    "C",
    "C",
    "C",
    "PF,C",
    "PF",
    "PG",
    "SG,SF",
    "SF,PF",
    "SG,SF",
    "SG,SF",
    "SG,SF,PF",
    "PG,SG"
)

$teams = @(
    "Sacramento Kings",
    "Detroit Pistons",
    "Boston Celtics",
    "New Orleans Pelicans",
    "Miami Heat",
    "Toronto Raptors",
    "Miami Heat",
    "Oklahoma City Thunder",
    "Los Angeles Lakers",
    "Minnesota Timberwolves",
    "Milwaukee Bucks",
    "Orlando Magic",
    "Brooklyn Nets",
    "Washington Wizards",
    "Charlotte Hornets",
    "New Orleans Pelicans",
    "Charlotte Hornets"
)

for ($i = 0; $i -lt $players.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $players[$i]
}

for ($i = 0; $i -lt $positions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $positions[$i]
}

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $teams[$i]
}
